$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.244.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.512.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.31%  "
$ws.Range("E7").Value = "  -0.57%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.505.80"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.194"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.20"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.584"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.33"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000276"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.077.67"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "610.48"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.533.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.422.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.42%  "
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.876"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.78%  "
$ws.Range("E23").Value = "  -19.39%  "
$ws.Range("E24").Value = "  -1.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.10%  "
$ws.Range("E26").Value = "  -3.84%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.38%  "
$ws.Range("E31").Value = "  -4.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "636.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +12.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.87"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.96%  "
$ws.Range("E35").Value = "  -4.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.60"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0995"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.71"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0471"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "56.73"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.48%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("E42").Value = "  +1.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₃0744"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.355.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.01%  "
$ws.Range("E45").Value = "  -5.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "32.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.60%  "
